{"js": "// Replace each two-digit multiplication expression's text with its\n// updated counterpart, per the unified diff. Each \"old\" text string\n// appears exactly once in the document, so a scoped, case-sensitive\n// search-and-replace for each pair unambiguously reproduces the edit.\nconst replacements = [\n  [\"93\u00d711=\", \"95\u00d723=\"],\n  [\"18\u00d797=\", \"80\u00d773=\"],\n  [\"21\u00d740=\", \"68\u00d775=\"],\n  [\"62\u00d773=\", \"69\u00d786=\"],\n  [\"38\u00d791=\", \"30\u00d773=\"],\n  [\"20\u00d767=\", \"18\u00d767=\"],\n  [\"14\u00d758=\", \"93\u00d747=\"],\n  [\"98\u00d781=\", \"52\u00d762=\"],\n  [\"41\u00d748=\", \"47\u00d751=\"],\n  [\"72\u00d743=\", \"12\u00d754=\"],\n  [\"66\u00d716=\", \"47\u00d747=\"],\n  [\"64\u00d780=\", \"97\u00d716=\"],\n  [\"96\u00d780=\", \"92\u00d717=\"],\n  [\"93\u00d789=\", \"85\u00d775=\"],\n  [\"67\u00d749=\", \"94\u00d725=\"],\n  [\"85\u00d729=\", \"97\u00d770=\"],\n  [\"24\u00d779=\", \"61\u00d717=\"],\n  [\"25\u00d787=\", \"30\u00d733=\"],\n  [\"55\u00d741=\", \"67\u00d713=\"],\n  [\"41\u00d774=\", \"21\u00d760=\"],\n  [\"68\u00d779=\", \"91\u00d718=\"],\n  [\"88\u00d760=\", \"75\u00d743=\"],\n  [\"15\u00d735=\", \"75\u00d767=\"],\n  [\"83\u00d784=\", \"89\u00d779=\"],\n  [\"96\u00d788=\", \"54\u00d739=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression's text with its\n# updated counterpart, per the unified diff. Each \"old\" text string\n# appears exactly once in the document, so a simple Find/Replace\n# (wdReplaceAll = 2, though only a single match ever exists) is\n# sufficient and unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"93\u00d711=\", \"95\u00d723=\"),\n    @(\"18\u00d797=\", \"80\u00d773=\"),\n    @(\"21\u00d740=\", \"68\u00d775=\"),\n    @(\"62\u00d773=\", \"69\u00d786=\"),\n    @(\"38\u00d791=\", \"30\u00d773=\"),\n    @(\"20\u00d767=\", \"18\u00d767=\"),\n    @(\"14\u00d758=\", \"93\u00d747=\"),\n    @(\"98\u00d781=\", \"52\u00d762=\"),\n    @(\"41\u00d748=\", \"47\u00d751=\"),\n    @(\"72\u00d743=\", \"12\u00d754=\"),\n    @(\"66\u00d716=\", \"47\u00d747=\"),\n    @(\"64\u00d780=\", \"97\u00d716=\"),\n    @(\"96\u00d780=\", \"92\u00d717=\"),\n    @(\"93\u00d789=\", \"85\u00d775=\"),\n    @(\"67\u00d749=\", \"94\u00d725=\"),\n    @(\"85\u00d729=\", \"97\u00d770=\"),\n    @(\"24\u00d779=\", \"61\u00d717=\"),\n    @(\"25\u00d787=\", \"30\u00d733=\"),\n    @(\"55\u00d741=\", \"67\u00d713=\"),\n    @(\"41\u00d774=\", \"21\u00d760=\"),\n    @(\"68\u00d779=\", \"91\u00d718=\"),\n    @(\"88\u00d760=\", \"75\u00d743=\"),\n    @(\"15\u00d735=\", \"75\u00d767=\"),\n    @(\"83\u00d784=\", \"89\u00d779=\"),\n    @(\"96\u00d788=\", \"54\u00d739=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
